$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = '''55.904.66'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.Value = '''  +8.75%  '
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.Value = '''3.217.64'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.Value = '''  +3.85%  '
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.Value = '''  +0.01%  '
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.Value = '''395.48'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.Value = '''  +1.93%  '
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.Value = '''110.37'
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.Value = '''  +6.45%  '
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.Value = '''0.551'
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.Value = '''  +2.34%  '
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.Value = '''  -0.03%  '
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.Value = '''0.616'
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.Value = '''  +4.83%  '
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.Value = '''39.07'
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.Value = '''  +5.39%  '
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.Value = '''0.0907'
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.Value = '''  +5.69%  '
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.Value = '''  +2.08%  '
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.Value = '''3.732.58'
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.Value = '''  +3.93%  '
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.Value = '''8.04'
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.Value = '''  +3.33%  '
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.Value = '''18.98'
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.Value = '''  +2.42%  '
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.Value = '''3.230.47'
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.Value = '''  +4.60%  '
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.Value = '''1.04'
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.Value = '''  +4.57%  '
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.Value = '''10.80'
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.Value = '''  +1.28%  '
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.Value = '''55.864.11'
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.Value = '''  +8.47%  '
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.Value = '''3.31'
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.Value = '''  +1.68%  '
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.Value = '''  +5.34%  '
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.Value = '''12.91'
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.Value = '''  +3.12%  '
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.Value = '''298.02'
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.Value = '''  +11.94%  '
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.Value = '''75.36'
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.Value = '''  +7.33%  '
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.Value = '''3.21'
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.Value = '''  +1.31%  '
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.Value = '''8.11'
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.Value = '''  +1.23%  '
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.Value = '''28.01'
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.Value = '''  +2.23%  '
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.Value = '''7.46'
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.Value = '''  +4.08%  '
$c.Style = 'Normal'
$c = $ws.Range('D29')
$c.Value = '''0.171'
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.Value = '''  +3.83%  '
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.Value = '''  +0.48%  '
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.Value = '''  +3.02%  '
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.Value = '''11.07'
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.Value = '''  +6.36%  '
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.Value = '''0.0489'
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.Value = '''  +3.51%  '
$c.Style = 'Normal'
$c = $ws.Range('D34')
$c.Value = '''35.94'
$c.Style = 'Normal'
$c = $ws.Range('E34')
$c.Value = '''  +0.58%  '
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.Value = '''  +2.61%  '
$c.Style = 'Normal'
$c = $ws.Range('D36')
$c.Value = '''51.32'
$c.Style = 'Normal'
$c = $ws.Range('B37')
$c.Value = '''LidoDAOToken'
$c.Style = 'Normal'
$c = $ws.Range('C37')
$c.Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c.Style = 'Normal'
$c = $ws.Range('D37')
$c.Value = '''3.52'
$c.Style = 'Normal'
$c = $ws.Range('E37')
$c.Value = '''  +4.32%  '
$c.Style = 'Normal'
$c = $ws.Range('B38')
$c.Value = '''FirstDigitalUSD'
$c.Style = 'Normal'
$c = $ws.Range('C38')
$c.Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.Value = '''1.00'
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.Value = '''  +0.06%  '
$c.Style = 'Normal'
$c = $ws.Range('B39')
$c.Value = '''Stacks'
$c.Style = 'Normal'
$c = $ws.Range('C39')
$c.Value = '''https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.Value = '''3.09'
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.Value = '''  +23.81%  '
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.Value = '''133.72'
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.Value = '''  +3.26%  '
$c.Style = 'Normal'
$c = $ws.Range('B41')
$c.Value = '''Celestia'
$c.Style = 'Normal'
$c = $ws.Range('C41')
$c.Value = '''https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.Value = '''17.29'
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.Value = '''  +4.43%  '
$c.Style = 'Normal'
$c = $ws.Range('B42')
$c.Value = '''ARBITRUM'
$c.Style = 'Normal'
$c = $ws.Range('C42')
$c.Value = '''https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.Value = '''1.91'
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.Value = '''  +2.87%  '
$c.Style = 'Normal'
$c = $ws.Range('D43')
$c.Value = '''3.97'
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.Value = '''  +3.86%  '
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.Value = '''  +2.84%  '
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.Value = '''  -3.29%  '
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.Value = '''22.15'
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.Value = '''  +0.27%  '
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.Value = '''2.18'
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.Value = '''  +52.39%  '
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.Value = '''  +1.64%  '
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.Value = '''  -1.60%  '
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.Value = '''2.126.24'
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.Value = '''  +2.39%  '
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.Value = '''0.0358'
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.Value = '''  +7.84%  '
$c.Style = 'Normal'
